$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly generated player-stat rows appended below the existing data.
# Columns: A=Player, B=Season Group, C=FG%, D=Lng
# Row 5-7   -> Daniel Carlson   (pale-yellow highlight, new style)
# Row 8-10  -> Harrison Butker  (existing green highlight, same as Ka'imi Fairbairn)
# Row 11-13 -> BrandonMcManus   (pale-yellow highlight, new style)

$data = @(
    @(5,  "Daniel Carlson",  "Group1",     86.8,                52.66666666666666),
    @(6,  "Daniel Carlson",  "Group2",     87.86666666666667,   55),
    @(7,  "Daniel Carlson",  "Difference", 1.066666666666677,   2.333333333333336),
    @(8,  "Harrison Butker", "Group1",     90.46666666666665,   56.66666666666666),
    @(9,  "Harrison Butker", "Group2",     84.43333333333334,   58.33333333333334),
    @(10, "Harrison Butker", "Difference", -6.033333333333317,  1.666666666666671),
    @(11, "BrandonMcManus",  "Group1",     83.86666666666666,   57.33333333333334),
    @(12, "BrandonMcManus",  "Group2",     84.69999999999999,   55.33333333333334),
    @(13, "BrandonMcManus",  "Difference", 0.8333333333333286,  -2)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# First apply the brand-new pale-yellow fill to row 5 (this mints the new style).
$ws.Range("A5:D5").Interior.Color = 12451839

# Propagate that exact style (same style index, no duplicate fill) to the
# other Daniel Carlson / BrandonMcManus rows via copy/paste-format.
$ws.Range("A5:D5").Copy()
$ws.Range("A6:D7").PasteSpecial(-4122)
$ws.Range("A11:D13").PasteSpecial(-4122)

# Harrison Butker rows reuse the pre-existing green highlight style
# (same as the original Ka'imi Fairbairn rows).
$ws.Range("A2:D2").Copy()
$ws.Range("A8:D10").PasteSpecial(-4122)

$excel.CutCopyMode = 0
